$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Friday hours for the week of row 7
$ws.Range("F7").Value = 6.25

# Update the selection to match the saved state
$ws.Range("G10").Select()
